$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()

$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.01475210851223867
$ws.Range("E2").Value = 0.0145633818779396
$ws.Range("D3").Value = 0.05054682952072559
$ws.Range("E3").Value = -0.01073454287079334
$ws.Range("D4").Value = 0.01423260854744168
$ws.Range("E4").Value = 0.02782309817485706
$ws.Range("D5").Value = 0.009659931403340961
$ws.Range("E5").Value = -0.01448717948717948
$ws.Range("D6").Value = 0.01561131606471981
$ws.Range("E6").Value = -0.01207012811867847
$ws.Range("D7").Value = 0.02009897342794369
$ws.Range("E7").Value = -0.02273385461917188
$ws.Range("D8").Value = 0.004567351797301449
$ws.Range("E8").Value = 0.01599126886210489
$ws.Range("D9").Value = 0.006569496425918264
$ws.Range("E9").Value = -0.006447234475738139
$ws.Range("D10").Value = 0.01421096658573997
$ws.Range("E10").Value = 0.009431557481519137
$ws.Range("D11").Value = 0.008113723152434398
$ws.Range("E11").Value = 0.0139662672670382
$ws.Range("D12").Value = 0.01548189775296992
$ws.Range("E12").Value = 0.01295896328293722
$ws.Range("D13").Value = 0.003182668745533443
$ws.Range("E13").Value = 0.02145045965270675
$ws.Range("D14").Value = 0.005810201048560784
$ws.Range("E14").Value = 0.004103165298944944
$ws.Range("D15").Value = 0.01453036354842223
$ws.Range("E15").Value = 0.01557189643452994
$ws.Range("D16").Value = 0.0108005463805816
$ws.Range("E16").Value = 0.01599587203302399
$ws.Range("D17").Value = 0.02098109235618916
$ws.Range("E17").Value = 0.001195298492595809
$ws.Range("D18").Value = 0.008599382396001085
$ws.Range("E18").Value = -0.006351126568879595
$ws.Range("D19").Value = 0.01675081643449083
$ws.Range("E19").Value = -0.001537824569704749
$ws.Range("D20").Value = 0.01244583085100961
$ws.Range("E20").Value = -0.008905915717199897
$ws.Range("D21").Value = 0.006770311538360795
$ws.Range("E21").Value = 0.005414551607444862
$ws.Range("D22").Value = 0.01513983710507213
$ws.Range("E22").Value = -0.01030691708657805
$ws.Range("D23").Value = 0.01922865076208499
$ws.Range("E23").Value = 0.003391009477436757
$ws.Range("D24").Value = 0.009961061188048953
$ws.Range("E24").Value = 0.02361636920858134
$ws.Range("D25").Value = 0.02032771564976383
$ws.Range("E25").Value = 0.02708696340883887
$ws.Range("D26").Value = 0.01390166301292146
$ws.Range("E26").Value = -0.004723820214208918
$ws.Range("D27").Value = 0.02173001569181549
$ws.Range("E27").Value = 0.01200551689549267
$ws.Range("D28").Value = 0.05498420570254235
$ws.Range("E28").Value = -0.01237682301931409
$ws.Range("D29").Value = 0.02047252173685814
$ws.Range("E29").Value = -0.01073304407398956
$ws.Range("D30").Value = 0.03037305354704321
$ws.Range("E30").Value = 0.004811416921508771
$ws.Range("D31").Value = 0.01524637500237628
$ws.Range("E31").Value = 0.004132541888947294
$ws.Range("D32").Value = 0.01327714225543174
$ws.Range("E32").Value = 0.005848471422241985
$ws.Range("D33").Value = 0.01833306366024832
$ws.Range("E33").Value = 0.02054961089494189
$ws.Range("D34").Value = 0.04421852176670485
$ws.Range("E34").Value = -0.007406598300221612
$ws.Range("D35").Value = 0.01083026924643804
$ws.Range("E35").Value = -0.003430531732418474
$ws.Range("D36").Value = 0.009943537081721099
$ws.Range("E36").Value = -0.02074978204010458
$ws.Range("D37").Value = 0.01102764765251592
$ws.Range("E37").Value = 0.01179195620130558
$ws.Range("D38").Value = 0.007393562881787889
$ws.Range("E38").Value = -0.001130653266331549
$ws.Range("D39").Value = 0.01210500865585584
$ws.Range("E39").Value = 0.01213130352045666
$ws.Range("D40").Value = 0.01736453169178899
$ws.Range("E40").Value = 0.0001890001890003656
$ws.Range("D41").Value = 0.01682202746727187
$ws.Range("E41").Value = 0.01126031612812994
$ws.Range("D42").Value = 0.03328310788229968
$ws.Range("E42").Value = -0.006976809086596147
$ws.Range("D43").Value = 0.01138528184366974
$ws.Range("E43").Value = 0.001446729358272281
$ws.Range("D44").Value = 0.02179419850527423
$ws.Range("E44").Value = -0.0009248249438498224
$ws.Range("D45").Value = 0.01283163984208984
$ws.Range("E45").Value = 0.01561380268844048
$ws.Range("D46").Value = 0.008659106779581356
$ws.Range("E46").Value = 0.01557522123893795
$ws.Range("D47").Value = 0.01324766708012411
$ws.Range("E47").Value = 0.01501829960876711
$ws.Range("D48").Value = 0.01033742697695605
$ws.Range("E48").Value = 0.01172569949862523
$ws.Range("D49").Value = 0.01628683011394446
$ws.Range("E49").Value = -0.0003421806029222374
$ws.Range("D50").Value = 0.008460427998122256
$ws.Range("E50").Value = 0.002923976608187218
$ws.Range("D51").Value = 0.01109350237717908
$ws.Range("E51").Value = 0.01368119630925846
$ws.Range("D52").Value = 0.008242584160449479
$ws.Range("E52").Value = -0.004462441120568594
$ws.Range("D53").Value = 0.009510945538235586
$ws.Range("E53").Value = 0.001184942120135091
$ws.Range("D54").Value = 0.1345419145971657
$ws.Range("E54").Value = 0.00009849305623954585
$ws.Range("D55").Value = 0.04392853805669301
$ws.Range("E55").Value = 0.0009585430146179519
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0.001287278743451026

$ws.Protect()
